$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fill in "Objetivos:" value (row 10) with the practical-experiments text
$ws.Range("B10:C10").Value = 'Apresentar aos alunos, na prática, experimentos envolvendo tópicos de Pirometalurgia; Solidificação / Fundição e Eletrometalurgia.'

# 2. Insert 3 new rows at row 13 for the additional professors, pushing the rest down
$ws.Rows("13:15").Insert()
# The insert carries column A's bold style down into A13:A15 as empty cells; drop it so
# those rows stay label-less (matches the target: only B/C populated on 13-15).
$ws.Range("A13:A15").Clear()

# 3. Fill the three professor rows (13-15), matching the B/C (normal + red) styling
#    used throughout the rest of the sheet.
$ws.Range("B13:C13").Value = '144651 - Antonio Fernando Sartori'
$ws.Range("B13").WrapText = $true
$ws.Range("B13").VerticalAlignment = -4160
$ws.Range("B13").Font.Bold = $false
$ws.Range("C13").WrapText = $true
$ws.Range("C13").VerticalAlignment = -4160
$ws.Range("C13").Font.Color = 255

$ws.Range("B14:C14").Value = '3577649 - Carlos Angelo Nunes'
$ws.Range("B14").WrapText = $true
$ws.Range("B14").VerticalAlignment = -4160
$ws.Range("B14").Font.Bold = $false
$ws.Range("C14").WrapText = $true
$ws.Range("C14").VerticalAlignment = -4160
$ws.Range("C14").Font.Color = 255

$ws.Range("B15:C15").Value = '5009972 - Gilberto Carvalho Coelho'
$ws.Range("B15").WrapText = $true
$ws.Range("B15").VerticalAlignment = -4160
$ws.Range("B15").Font.Bold = $false
$ws.Range("C15").WrapText = $true
$ws.Range("C15").VerticalAlignment = -4160
$ws.Range("C15").Font.Color = 255

# 4. "Programa resumido:" value (now row 16) -- cell already styled, just set the text
$ws.Range("B16:C16").Value = 'Aulas práticas: pirometalurgia; eletrometalurgia e solidificação / fundição de metais e ligas.'

# 5. "Programa:" value (now row 18)
$ws.Range("B18:C18").Value = 'Eletrometalurgia:1) Eletrodeposição de cobre e/ou níquel;2) Anodização.PirometalurgiaOs experimentos poderão ser alterados a cada oferecimento da disciplina, mas envolverão normalmente os seguintes assuntos:1) Caracterização química e microestrutural de minérios; 2) Caracterização química e microestrutural de Ferro-Ligas;3) Secagem; calcinação de carbonatos e/ou hidróxidos; ustulação de sulfetos; aglomeração de minérios;4) Redução carbotérmica, por hidrogênio ou metalotérmica de óxidos; 5) Decomposição térmica sob vácuo de óxidos; 6) Refino de metais e/ou ligas.Solidificação / Fundição:Verificação experimental do efeito de algumas variáveis básicas de fundição tais como tipo de molde (cerâmicos, metálicos, ...), refrigeração ou pré-aquecimento do molde, temperatura de vazamento, adição de agentes inoculantes, agitação mecânica, dentre outras, nas micro e macroestruturas de metais e ligas metálicas. A turma será dividida em grupos sendo que cada um dos grupos deverá investigar com mais detalhe uma das variáveis (ou combinação de variáveis) do processamento por fundição.'

# 6. "Método:" value (now row 21)
$ws.Range("B21:C21").Value = 'Experimentos desenvolvidos em laboratórios, realização de relatórios para cada experimento e apresentação dos resultados obtidos.'

# 7. "Critério:" value (now row 22)
$ws.Range("B22:C22").Value = 'A avaliação será feita através das notas dos relatórios escritos e das apresentações orais dos trabalhos desenvolvidos em cada módulo da disciplina.'

# 8. "Norma de recuperação:" value (now row 23)
$ws.Range("B23:C23").Value = 'Devido às características da disciplina não será oferecida recuperação'

# 9. "Bibliografia:" value (now row 24)
$ws.Range("B24:C24").Value = 'DENARO, A.R. Fundamentos de Eletroquímica. São Paulo: Editora Edgard Blucher, 1974. PLETCHER, D.; WALSH, F.C. Industrial Electrochemistry, 2ª Ed. Springer, 1990.KUHN, A.T. Industrial Electrochemical Processes. Elsevier Pub., 1971.ROSENQVIST, T. Principles of Extrative Metallurgy, McGraw-Hill, 1983.HABASHI, F. Extractive Metallurgy, Gordon and Breach Science Publishers, 1986.GARCIA, A. Solidificação: Fundamentos e Aplicações, Editora da Unicamp, 2001. FLEMINGS, M.C. Solidification Processing, McGraw-Hill, 1974.'

# 10. Column A had shared a width group with column B ("1,2"); the edit splits it into
#     its own single-column width entry ("1,1") while keeping the same 30.7109375 width.
$ws.Columns("A").ColumnWidth = 30.7109375
